$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.454.33"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.724.73"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").Value = "'0.9967"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'242.85"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").Value = "'0.9975"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.4888"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").Value = "'0.2610"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "'0.06196"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "1.725.09"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'0.06985"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "'15.59"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "'4.520"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "'0.6012"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "'77.23"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "'0.9972"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "26.435.85"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "'0.9968"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'0.000007162"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "'11.32"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").Value = "1.939.22"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "'4.463"
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").Value = "'8.517"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").Value = "'5.106"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "'137.63"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "'15.27"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "'1.410"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'106.59"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").Value = "'1.743"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("D30").Value = "'3.913"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("D31").Value = "'0.08021"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "'3.648"
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("D33").Value = "'0.04492"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "'0.9964"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.600"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.6246"
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.9274"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'1.984"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.385"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'0.9971"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01480"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'99.83"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.378"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.3843"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'6.908"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1162"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05363"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'30.39"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.711"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'51.23"
$ws.Range("E51").Value = "  -0.36%  "
